# Auto-generated Excel COM-interop script applying scheduled market-data refresh
# to the Tonberry_Profits workbook leve-profit tables (columns H-N) across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H18").Value = 14688.048
$ws.Range("I18").Value = 10831.111
$ws.Range("J18").Value = 17580.75
$ws.Range("K18").Value = 10831.111
$ws.Range("L18").Value = 17580.75
$ws.Range("M18").Value = -10547.111
$ws.Range("N18").Value = -18148.75
$ws.Range("H28").Value = 3804527.5
$ws.Range("I28").Value = 4750160
$ws.Range("J28").Value = 21999
$ws.Range("K28").Value = 4750160
$ws.Range("L28").Value = 21999
$ws.Range("M28").Value = -4749675
$ws.Range("N28").Value = -22969
$ws.Range("H31").Value = 92.666664
$ws.Range("I31").Value = 136.5
$ws.Range("J31").Value = 5
$ws.Range("K31").Value = 409.5
$ws.Range("L31").Value = 15
$ws.Range("M31").Value = -179.5
$ws.Range("N31").Value = -475
$ws.Range("H42").Value = 216
$ws.Range("I42").Value = 65.666664
$ws.Range("J42").Value = 366.33334
$ws.Range("K42").Value = 196.999992
$ws.Range("L42").Value = 1099.00002
$ws.Range("M42").Value = 33.00000800000001
$ws.Range("N42").Value = -1559.00002
$ws.Range("H86").Value = 2372.25
$ws.Range("I86").Value = 2166.3333
$ws.Range("K86").Value = 2166.3333
$ws.Range("M86").Value = -1043.3333
$ws.Range("H89").Value = 2372.25
$ws.Range("I89").Value = 2166.3333
$ws.Range("K89").Value = 10831.6665
$ws.Range("M89").Value = -5215.666499999999
$ws.Range("H106").Value = 2495.1667
$ws.Range("I106").Value = 2495.1667
$ws.Range("K106").Value = 2495.1667
$ws.Range("M106").Value = -1864.1667
$ws.Range("H113").Value = 16506.25
$ws.Range("I113").Value = 18364.285
$ws.Range("K113").Value = 18364.285
$ws.Range("M113").Value = -15110.285
$ws.Range("H132").Value = 966.1591
$ws.Range("I132").Value = 975.8333
$ws.Range("K132").Value = 2927.4999
$ws.Range("M132").Value = -397.4998999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5132.425
$ws.Range("I32").Value = 3579
$ws.Range("J32").Value = 9227.817999999999
$ws.Range("K32").Value = 3579
$ws.Range("L32").Value = 9227.817999999999
$ws.Range("M32").Value = -3292
$ws.Range("N32").Value = -9801.817999999999
$ws.Range("H63").Value = 3870
$ws.Range("I63").Value = 3870
$ws.Range("K63").Value = 3870
$ws.Range("M63").Value = -3184
$ws.Range("H66").Value = 3870
$ws.Range("I66").Value = 3870
$ws.Range("K66").Value = 19350
$ws.Range("M66").Value = -15918
$ws.Range("H74").Value = 418.4375
$ws.Range("I74").Value = 418.4375
$ws.Range("K74").Value = 418.4375
$ws.Range("M74").Value = 455.5625
$ws.Range("H77").Value = 418.4375
$ws.Range("I77").Value = 418.4375
$ws.Range("K77").Value = 2092.1875
$ws.Range("M77").Value = 2275.8125
$ws.Range("H101").Value = 49980
$ws.Range("J101").Value = 49980
$ws.Range("L101").Value = 49980
$ws.Range("N101").Value = -56470
$ws.Range("H132").Value = 1781.4286
$ws.Range("I132").Value = 1465.5294
$ws.Range("J132").Value = 3124
$ws.Range("K132").Value = 4396.5882
$ws.Range("L132").Value = 9372
$ws.Range("M132").Value = -1866.5882
$ws.Range("N132").Value = -14432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 29665.334
$ws.Range("J76").Value = 29665.334
$ws.Range("L76").Value = 29665.334
$ws.Range("N76").Value = -30295.334
$ws.Range("H79").Value = 29665.334
$ws.Range("J79").Value = 29665.334
$ws.Range("L79").Value = 29665.334
$ws.Range("N79").Value = -31849.334
$ws.Range("H82").Value = 34999.668
$ws.Range("H85").Value = 34999.668
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H107").Value = 2749.6667
$ws.Range("I107").Value = 2749.6667
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2749.6667
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -829.6667000000002
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2977
$ws.Range("I31").Value = 1368.375
$ws.Range("J31").Value = 7266.6665
$ws.Range("K31").Value = 1368.375
$ws.Range("L31").Value = 7266.6665
$ws.Range("M31").Value = -1073.375
$ws.Range("N31").Value = -7856.6665
$ws.Range("H34").Value = 2977
$ws.Range("I34").Value = 1368.375
$ws.Range("J34").Value = 7266.6665
$ws.Range("K34").Value = 1368.375
$ws.Range("L34").Value = 7266.6665
$ws.Range("M34").Value = -1166.375
$ws.Range("N34").Value = -7670.6665
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H58").Value = 1153.0385
$ws.Range("I58").Value = 860.55554
$ws.Range("K58").Value = 860.55554
$ws.Range("M58").Value = -657.55554
$ws.Range("H132").Value = 1758.6428
$ws.Range("I132").Value = 1051.25
$ws.Range("J132").Value = 6003
$ws.Range("K132").Value = 3153.75
$ws.Range("L132").Value = 18009
$ws.Range("M132").Value = -623.75
$ws.Range("N132").Value = -23069
$ws.Range("H134").Value = 868.875
$ws.Range("I134").Value = 723.46155
$ws.Range("J134").Value = 1499
$ws.Range("K134").Value = 2170.38465
$ws.Range("L134").Value = 4497
$ws.Range("M134").Value = 364.61535
$ws.Range("N134").Value = -9567
$ws.Range("H136").Value = 1153.0385
$ws.Range("I136").Value = 860.55554
$ws.Range("K136").Value = 2581.66662
$ws.Range("M136").Value = -31.66661999999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1.6666666
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 6
$ws.Range("M16").Value = 170
$ws.Range("N16").Value = -352
$ws.Range("H23").Value = 135.75
$ws.Range("I23").Value = 49
$ws.Range("K23").Value = 147
$ws.Range("M23").Value = 88
$ws.Range("H131").Value = 5690533.5
$ws.Range("J131").Value = 9104.487999999999
$ws.Range("L131").Value = 27313.464
$ws.Range("N131").Value = -37393.464

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 600
$ws.Range("I23").Value = 600
$ws.Range("K23").Value = 600
$ws.Range("M23").Value = -370
$ws.Range("H46").Value = 1606.6
$ws.Range("I46").Value = 1188.6666
$ws.Range("K46").Value = 1188.6666
$ws.Range("M46").Value = -1000.6666
$ws.Range("H82").Value = 5616
$ws.Range("J82").Value = 6520
$ws.Range("L82").Value = 6520
$ws.Range("N82").Value = -7242
$ws.Range("H85").Value = 5616
$ws.Range("J85").Value = 6520
$ws.Range("L85").Value = 6520
$ws.Range("N85").Value = -9016

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 593.46155
$ws.Range("J107").Value = 692
$ws.Range("L107").Value = 2076
$ws.Range("N107").Value = -5916
$ws.Range("H113").Value = 1107.7142
$ws.Range("I113").Value = 1038.5
$ws.Range("K113").Value = 3115.5
$ws.Range("M113").Value = -945.5
$ws.Range("H126").Value = 5129.3335
$ws.Range("I126").Value = 5026.684
$ws.Range("J126").Value = 5519.4
$ws.Range("K126").Value = 15080.052
$ws.Range("L126").Value = 16558.2
$ws.Range("M126").Value = -12610.052
$ws.Range("N126").Value = -21498.2
$ws.Range("H136").Value = 4443.278
$ws.Range("I136").Value = 4998.2
$ws.Range("K136").Value = 14994.6
$ws.Range("M136").Value = -12444.6
